$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.040226340180912
$ws.Range("D2").Value = 1.047336845981845
$ws.Range("E2").Value = 1.03863965568597
$ws.Range("F2").Value = 1.055638296016008
$ws.Range("I2").Value = 1.036321019683814
$ws.Range("J2").Value = 1.045314091273764
$ws.Range("K2").Value = 1.050099766748066
$ws.Range("L2").Value = 1.04142711491686
$ws.Range("M2").Value = 1.058378203606451
$ws.Range("N2").Value = 1.018955348806655
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041351736615342
$ws.Range("D3").Value = 1.048351073446975
$ws.Range("E3").Value = 1.039601712462545
$ws.Range("F3").Value = 1.05674010434856
$ws.Range("I3").Value = 1.036506607541424
$ws.Range("J3").Value = 1.046084139735947
$ws.Range("K3").Value = 1.050925293048035
$ws.Range("L3").Value = 1.042198825033247
$ws.Range("M3").Value = 1.059292757763106
$ws.Range("N3").Value = 1.019216001933543
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042079954943642
$ws.Range("D4").Value = 1.049007641102278
$ws.Range("E4").Value = 1.040224574054876
$ws.Range("F4").Value = 1.057453425480667
$ws.Range("I4").Value = 1.036625204940288
$ws.Range("J4").Value = 1.046581879080664
$ws.Range("K4").Value = 1.051459149837395
$ws.Range("L4").Value = 1.042697901974708
$ws.Range("M4").Value = 1.059884320470822
$ws.Range("N4").Value = 1.019384351613516
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042386101411293
$ws.Range("D5").Value = 1.049283732734614
$ws.Range("E5").Value = 1.040486507960088
$ws.Range("F5").Value = 1.057753396188461
$ws.Range("I5").Value = 1.036674706369538
$ws.Range("J5").Value = 1.046791001087491
$ws.Range("K5").Value = 1.051683508160625
$ws.Range("L5").Value = 1.042907649041724
$ws.Range("M5").Value = 1.060132962059257
$ws.Range("N5").Value = 1.019455051473494
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042437504975233
$ws.Range("D6").Value = 1.049330093875346
$ws.Range("E6").Value = 1.040530492674196
$ws.Range("F6").Value = 1.057803767901466
$ws.Range("I6").Value = 1.036682996943949
$ws.Range("J6").Value = 1.046826106112434
$ws.Range("K6").Value = 1.051721174510573
$ws.Range("L6").Value = 1.042942862693374
$ws.Range("M6").Value = 1.0601747070662
$ws.Range("N6").Value = 1.019466917929856
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04208404567315
$ws.Range("D7").Value = 1.049011329974755
$ws.Range("E7").Value = 1.040228073702916
$ws.Range("F7").Value = 1.057457433348155
$ws.Range("I7").Value = 1.036625867783403
$ws.Range("J7").Value = 1.046584673880458
$ws.Range("K7").Value = 1.051462148019167
$ws.Range("L7").Value = 1.042700704880187
$ws.Range("M7").Value = 1.059887643033646
$ws.Range("N7").Value = 1.019385296600701
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.040606671203577
$ws.Range("D8").Value = 1.047679547918218
$ws.Range("E8").Value = 1.038964715487316
$ws.Range("F8").Value = 1.056010579227328
$ws.Range("I8").Value = 1.036384048483159
$ws.Range("J8").Value = 1.045574443287322
$ws.Range("K8").Value = 1.05037882227814
$ws.Range("L8").Value = 1.041687974113408
$ws.Range("M8").Value = 1.058687326037041
$ws.Range("N8").Value = 1.019043501973686
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.038003405212314
$ws.Range("D9").Value = 1.045335031896258
$ws.Range("E9").Value = 1.036741173283882
$ws.Range("F9").Value = 1.053463921589116
$ws.Range("I9").Value = 1.035946522867552
$ws.Range("J9").Value = 1.043790196599426
$ws.Range("K9").Value = 1.048467455316499
$ws.Range("L9").Value = 1.039901335153181
$ws.Range("M9").Value = 1.056570562825795
$ws.Range("N9").Value = 1.018438841630936
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.036267884919996
$ws.Range("D10").Value = 1.043773529903556
$ws.Range("E10").Value = 1.035260598102259
$ws.Range("F10").Value = 1.051768076325791
$ws.Range("I10").Value = 1.035647174234391
$ws.Range("J10").Value = 1.042597937984961
$ws.Range("K10").Value = 1.047191582562689
$ws.Range("L10").Value = 1.038708840554259
$ws.Range("M10").Value = 1.055158267500551
$ws.Range("N10").Value = 1.018034140169542
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035516371705961
$ws.Range("D11").Value = 1.043097737640106
$ws.Range("E11").Value = 1.034619915486751
$ws.Range("F11").Value = 1.051034209713021
$ws.Range("I11").Value = 1.035515735844744
$ws.Range("J11").Value = 1.042081018946255
$ws.Range("K11").Value = 1.046638726609455
$ws.Range("L11").Value = 1.038192142796543
$ws.Range("M11").Value = 1.054546459353112
$ws.Range("N11").Value = 1.017858521850015
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035237221913917
$ws.Range("D12").Value = 1.042846770322877
$ws.Range("E12").Value = 1.034381999784564
$ws.Range("F12").Value = 1.050761685476959
$ws.Range("I12").Value = 1.035466640582665
$ws.Range("J12").Value = 1.04188891219726
$ws.Range("K12").Value = 1.046433311660152
$ws.Range("L12").Value = 1.038000166942029
$ws.Range("M12").Value = 1.054319164884388
$ws.Range("N12").Value = 1.017793232275295
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035297100705149
$ws.Range("D13").Value = 1.042900601315608
$ws.Range("E13").Value = 1.034433030691906
$ws.Range("F13").Value = 1.050820139848598
$ws.Range("I13").Value = 1.035477184045378
$ws.Range("J13").Value = 1.041930124278228
$ws.Range("K13").Value = 1.046477376568092
$ws.Range("L13").Value = 1.038041348740202
$ws.Range("M13").Value = 1.054367922211706
$ws.Range("N13").Value = 1.017807239695662
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035493297191602
$ws.Range("D14").Value = 1.043076991544083
$ws.Range("E14").Value = 1.034600248029782
$ws.Range("F14").Value = 1.051011681420656
$ws.Range("I14").Value = 1.035511683188035
$ws.Range("J14").Value = 1.042065141384586
$ws.Range("K14").Value = 1.046621748164138
$ws.Range("L14").Value = 1.038176275060568
$ws.Range("M14").Value = 1.054527671973886
$ws.Range("N14").Value = 1.017853126158654
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035614179716548
$ws.Range("D15").Value = 1.04318567826876
$ws.Range("E15").Value = 1.034703284406678
$ws.Range("F15").Value = 1.051129705303398
$ws.Range("I15").Value = 1.035532903048185
$ws.Range("J15").Value = 1.042148316617739
$ws.Range("K15").Value = 1.046710692357172
$ws.Range("L15").Value = 1.038259400813921
$ws.Range("M15").Value = 1.054626093540315
$ws.Range("N15").Value = 1.017881390754114
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036317759830868
$ws.Range("D16").Value = 1.043818387323198
$ws.Range("E16").Value = 1.035303126850867
$ws.Range("F16").Value = 1.051816789985902
$ws.Range("I16").Value = 1.035655859055556
$ws.Range("J16").Value = 1.042632230195301
$ws.Range("K16").Value = 1.047228265447795
$ws.Range("L16").Value = 1.038743124902811
$ws.Range("M16").Value = 1.055198865367566
$ws.Range("N16").Value = 1.018045787371804
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036759090524655
$ws.Range("D17").Value = 1.044215362046199
$ws.Range("E17").Value = 1.035679503470927
$ws.Range("F17").Value = 1.052247899035781
$ws.Range("I17").Value = 1.035732499265801
$ws.Range("J17").Value = 1.042935598585339
$ws.Range("K17").Value = 1.04755281947755
$ws.Range("L17").Value = 1.039046461105621
$ws.Range("M17").Value = 1.055558076126776
$ws.Range("N17").Value = 1.0181488072893
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037016509253166
$ws.Range("D18").Value = 1.04444694432363
$ws.Range("E18").Value = 1.035899077756098
$ws.Range("F18").Value = 1.052499400672343
$ws.Range("I18").Value = 1.035777026748445
$ws.Range("J18").Value = 1.043112484270916
$ws.Range("K18").Value = 1.047742088298991
$ws.Range("L18").Value = 1.039223359219539
$ws.Range("M18").Value = 1.05576757118151
$ws.Range("N18").Value = 1.018208860410118
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037104282081691
$ws.Range("D19").Value = 1.044525913611011
$ws.Range("E19").Value = 1.035973953732175
$ws.Range("F19").Value = 1.052585163541703
$ws.Range("I19").Value = 1.035792179708737
$ws.Range("J19").Value = 1.043172786880294
$ws.Range("K19").Value = 1.047806617631921
$ws.Range("L19").Value = 1.039283671355992
$ws.Range("M19").Value = 1.055838999097188
$ws.Range("N19").Value = 1.018229330756232
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.036711740132359
$ws.Range("D20").Value = 1.044172766938239
$ws.Range("E20").Value = 1.035639117667756
$ws.Range("F20").Value = 1.052201640667828
$ws.Range("I20").Value = 1.035724294644603
$ws.Range("J20").Value = 1.042903056646715
$ws.Range("K20").Value = 1.047518001840063
$ws.Range("L20").Value = 1.039013919392207
$ws.Range("M20").Value = 1.05551953895581
$ws.Range("N20").Value = 1.018137758018625
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035435522344556
$ws.Range("D21").Value = 1.04302504758484
$ws.Range("E21").Value = 1.034551004971939
$ws.Range("F21").Value = 1.050955275369591
$ws.Range("I21").Value = 1.035501531591739
$ws.Range("J21").Value = 1.042025384975463
$ws.Range("K21").Value = 1.046579235976377
$ws.Range("L21").Value = 1.03813654404044
$ws.Range("M21").Value = 1.054480630804213
$ws.Range("N21").Value = 1.01783961531582
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034633088013224
$ws.Range("D22").Value = 1.042303731679131
$ws.Range("E22").Value = 1.033867225965193
$ws.Range("F22").Value = 1.05017202056298
$ws.Range("I22").Value = 1.035359891020273
$ws.Range("J22").Value = 1.041472979197739
$ws.Range("K22").Value = 1.045988651433215
$ws.Range("L22").Value = 1.03758460639004
$ws.Range("M22").Value = 1.0538271866704
$ws.Range("N22").Value = 1.017651830718939
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03505847642334
$ws.Range("D23").Value = 1.042686086484542
$ws.Range("E23").Value = 1.034229676032082
$ws.Range("F23").Value = 1.050587202455592
$ws.Range("I23").Value = 1.035435127186136
$ws.Range("J23").Value = 1.041765874934736
$ws.Range("K23").Value = 1.046301764308371
$ws.Range("L23").Value = 1.037877227215974
$ws.Range("M23").Value = 1.054173612771067
$ws.Range("N23").Value = 1.017751410220656
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.036733135742278
$ws.Range("D24").Value = 1.044192013729903
$ws.Range("E24").Value = 1.035657366149533
$ws.Range("F24").Value = 1.052222542699973
$ws.Range("I24").Value = 1.035728002501828
$ws.Range("J24").Value = 1.042917761145468
$ws.Range("K24").Value = 1.047533734549118
$ws.Range("L24").Value = 1.039028623692811
$ws.Range("M24").Value = 1.055536952326236
$ws.Range("N24").Value = 1.018142750821636
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.038676410076927
$ws.Range("D25").Value = 1.045940879059485
$ws.Range("E25").Value = 1.037315697407887
$ws.Range("F25").Value = 1.05412195314032
$ws.Range("I25").Value = 1.036060985149884
$ws.Range("J25").Value = 1.044251953007518
$ws.Range("K25").Value = 1.04896187581291
$ws.Range("L25").Value = 1.040363470891795
$ws.Range("M25").Value = 1.057117993526121
$ws.Range("N25").Value = 1.018595442013909
